$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.871.15"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "3.119.69"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.38"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.26%  "
$ws.Range("E10").Value = "  -1.40%  "
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000249"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.03%  "
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("D15").Value = "3.633.73"
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").Value = "66.846.82"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").Value = "3.118.18"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "474.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.67%  "
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "83.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.62%  "
$ws.Range("E25").Value = "  -3.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.11%  "
$ws.Range("E29").Value = "  -1.45%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("E33").Value = "  -6.48%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("E36").Value = "  -3.31%  "
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("E39").Value = "  -1.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.312"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("E41").Value = "  +1.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("D43").Value = "2.814.34"
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "382.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0353"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "135.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("E51").Value = "  -0.86%  "
